$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = -6
$ws.Range("F2").Value = 1593.07
$ws.Range("G2").Value = 19.71
$ws.Range("V2").Value = [double]"1.255437614258792e-08"
$ws.Range("W2").Value = 0.0001255437614234308
$ws.Range("D3").Value = 99380
$ws.Range("E3").Value = 569169
$ws.Range("F3").Value = 1593.33
$ws.Range("G3").Value = 19.73
$ws.Range("P3").Value = 0.0009940168196829501
$ws.Range("Q3").Value = [double]"6.321253794347984e-08"
$ws.Range("R3").Value = 1.988966074789402
$ws.Range("S3").Value = 5.68601177492092
$ws.Range("T3").Value = 0.03371896508567265
$ws.Range("U3").Value = 2.08343460112531
$ws.Range("V3").Value = [double]"1.958623812165344e-08"
$ws.Range("W3").Value = 0.000178849360134071
$ws.Range("D4").Value = 198010
$ws.Range("E4").Value = 838356
$ws.Range("F4").Value = 1593.46
$ws.Range("G4").Value = 19.71
$ws.Range("P4").Value = 0.001980521284213081
$ws.Range("Q4").Value = [double]"9.853466689660475e-08"
$ws.Range("R4").Value = 2.000304638636801
$ws.Range("S4").Value = 8.374497161212975
$ws.Range("T4").Value = 0.04959054289304366
$ws.Range("U4").Value = 2.084499826805902
$ws.Range("V4").Value = [double]"3.156950443102184e-08"
$ws.Range("W4").Value = 0.0002622134371218887
$ws.Range("D5").Value = 296879
$ws.Range("E5").Value = 1066481
$ws.Range("F5").Value = 1593.67
$ws.Range("G5").Value = 19.71
$ws.Range("P5").Value = 0.002969416642769562
$ws.Range("Q5").Value = [double]"1.382915076109563e-07"
$ws.Range("R5").Value = 2.010169515041158
$ws.Range("S5").Value = 10.65326474493345
$ws.Range("T5").Value = 0.0630539012780922
$ws.Range("U5").Value = 2.084843858322008
$ws.Range("V5").Value = [double]"4.471547759106351e-08"
$ws.Range("W5").Value = 0.0003552745963425761
$ws.Range("D6").Value = 394832
$ws.Range("E6").Value = 1265232
$ws.Range("F6").Value = 1593.74
$ws.Range("G6").Value = 19.7
$ws.Range("P6").Value = 0.003949149308890628
$ws.Range("Q6").Value = [double]"1.790148686677638e-07"
$ws.Range("R6").Value = 2.015710172132742
$ws.Range("S6").Value = 12.63810584218668
$ws.Range("T6").Value = 0.07478944908286286
$ws.Range("U6").Value = 2.084988667016243
$ws.Range("V6").Value = [double]"5.80893985225214e-08"
$ws.Range("W6").Value = 0.0004506520815701664
$ws.Range("D7").Value = 493210
$ws.Range("E7").Value = 1442337
$ws.Range("F7").Value = 1592.85
$ws.Range("G7").Value = 19.71
$ws.Range("P7").Value = 0.00493313443743475
$ws.Range("Q7").Value = [double]"2.204442596719743e-07"
$ws.Range("R7").Value = 2.019069475039879
$ws.Range("S7").Value = 14.40773474579633
$ws.Range("T7").Value = 0.08524795306321042
$ws.Range("U7").Value = 2.08507527367559
$ws.Range("V7").Value = [double]"7.166292781028802e-08"
$ws.Range("W7").Value = 0.0005476367680598177
$ws.Range("D8").Value = 591464
$ws.Range("E8").Value = 1611164
$ws.Range("F8").Value = 1592.97
$ws.Range("G8").Value = 19.7
$ws.Range("P8").Value = 0.005915877335408122
$ws.Range("Q8").Value = [double]"2.621421023954124e-07"
$ws.Range("R8").Value = 2.021207037917273
$ws.Range("S8").Value = 16.09352308802535
$ws.Range("T8").Value = 0.09522098677452703
$ws.Range("U8").Value = 2.085115084871771
$ws.Range("V8").Value = [double]"8.530441989505636e-08"
$ws.Range("W8").Value = 0.0006454674611146695
$ws.Range("D9").Value = 689154
$ws.Range("E9").Value = 1775693
$ws.Range("F9").Value = 1593.08
$ws.Range("G9").Value = 19.7
$ws.Range("P9").Value = 0.006892980049668741
$ws.Range("Q9").Value = [double]"3.03798234077337e-07"
$ws.Range("R9").Value = 2.022640790532277
$ws.Range("S9").Value = 17.73695689594307
$ws.Range("T9").Value = 0.1049402095214556
$ws.Range("U9").Value = 2.085146043611084
$ws.Range("V9").Value = [double]"9.892234889453298e-08"
$ws.Range("W9").Value = 0.0007434423086112874
$ws.Range("D10").Value = 786400
$ws.Range("E10").Value = 1927498
$ws.Range("F10").Value = 1593.25
$ws.Range("G10").Value = 19.69
$ws.Range("P10").Value = 0.007865640269459415
$ws.Range("Q10").Value = [double]"3.453360394319712e-07"
$ws.Range("R10").Value = 2.023693546726142
$ws.Range("S10").Value = 19.25252319952976
$ws.Range("T10").Value = 0.1139101160682761
$ws.Range("U10").Value = 2.085155411240929
$ws.Range("V10").Value = [double]"1.124964921201696e-07"
$ws.Range("W10").Value = 0.0008409970778318691
$ws.Range("D11").Value = 883362
$ws.Range("E11").Value = 2073610
$ws.Range("F11").Value = 1593.16
$ws.Range("G11").Value = 19.7
$ws.Range("P11").Value = 0.00883546304619987
$ws.Range("Q11").Value = [double]"3.868176692643267e-07"
$ws.Range("R11").Value = 2.024487040268303
$ws.Range("S11").Value = 20.71276583997981
$ws.Range("T11").Value = 0.1225417546724499
$ws.Range("U11").Value = 2.085173993463791
$ws.Range("V11").Value = [double]"1.260500337935669e-07"
$ws.Range("W11").Value = 0.000938461759000262
$ws.Range("D12").Value = 980451
$ws.Range("E12").Value = 2222585
$ws.Range("F12").Value = 1592.88
$ws.Range("G12").Value = 19.86
$ws.Range("P12").Value = 0.009806585894470006
$ws.Range("Q12").Value = [double]"4.284138567750588e-07"
$ws.Range("R12").Value = 2.025100348748243
$ws.Range("S12").Value = 22.2150606926055
$ws.Range("T12").Value = 0.1313269406267873
$ws.Range("U12").Value = 2.085274685471221
$ws.Range("V12").Value = [double]"1.396457138224739e-07"
$ws.Range("W12").Value = 0.001036482195090958
$ws.Range("D13").Value = 99118
$ws.Range("E13").Value = 1166562
$ws.Range("F13").Value = 1594.81
$ws.Range("G13").Value = 19.77
$ws.Range("P13").Value = [double]"9.914170718239515e-05"
$ws.Range("Q13").Value = [double]"7.21351138133991e-09"
$ws.Range("R13").Value = 1.997657557397754
$ws.Range("S13").Value = 1.165597564479266
$ws.Range("T13").Value = 0.006895113510952668
$ws.Range("U13").Value = 2.085048485178118
$ws.Range("V13").Value = [double]"2.204574624214287e-09"
$ws.Range("W13").Value = [double]"2.285437881499195e-05"
$ws.Range("D14").Value = 198956
$ws.Range("E14").Value = 1739697
$ws.Range("F14").Value = 1594.95
$ws.Range("G14").Value = 19.76
$ws.Range("P14").Value = 0.0001990005243383183
$ws.Range("Q14").Value = [double]"1.11644392449385e-08"
$ws.Range("R14").Value = 2.002804983453184
$ws.Range("S14").Value = 1.738176858421052
$ws.Range("T14").Value = 0.01027892336698842
$ws.Range("U14").Value = 2.085310073182456
$ws.Range("V14").Value = [double]"3.515831111854671e-09"
$ws.Range("W14").Value = [double]"3.38604964230146e-05"
$ws.Range("D15").Value = 298184
$ws.Range("E15").Value = 2207103
$ws.Range("F15").Value = 1595.08
$ws.Range("G15").Value = 19.75
$ws.Range("P15").Value = 0.0002982491743606949
$ws.Range("Q15").Value = [double]"1.532759379682069e-08"
$ws.Range("R15").Value = 2.006586444293072
$ws.Range("S15").Value = 2.20507943285509
$ws.Range("T15").Value = 0.01303923098466725
$ws.Range("U15").Value = 2.085383461903509
$ws.Range("V15").Value = [double]"4.886582170902018e-09"
$ws.Range("W15").Value = [double]"4.473974635527407e-05"
$ws.Range("D16").Value = 396438
$ws.Range("E16").Value = 2792679
$ws.Range("F16").Value = 1595.22
$ws.Range("G16").Value = 19.74
$ws.Range("P16").Value = 0.0003965235817915425
$ws.Range("Q16").Value = [double]"1.980695884882506e-08"
$ws.Range("R16").Value = 2.009449846405277
$ws.Range("S16").Value = 2.790000203016651
$ws.Range("T16").Value = 0.01649771681718106
$ws.Range("U16").Value = 2.085427176881368
$ws.Range("V16").Value = [double]"6.339493352381419e-09"
$ws.Range("W16").Value = [double]"5.708173665713881e-05"
$ws.Range("D17").Value = 494367
$ws.Range("E17").Value = 3176412
$ws.Range("F17").Value = 1.59532
$ws.Range("G17").Value = 19.71
$ws.Range("P17").Value = 0.0004944726844883379
$ws.Range("Q17").Value = [double]"3.643750200513162e-07"
$ws.Range("R17").Value = 2.224972110923106
$ws.Range("S17").Value = 3.172980485914009
$ws.Range("T17").Value = 0.01899651737571835
$ws.Range("U17").Value = 2.077405926595824
$ws.Range("V17").Value = [double]"1.253690113276088e-06"
$ws.Range("W17").Value = 0.01253684615535567
$ws.Range("D18").Value = 593385
$ws.Range("E18").Value = 3538125
$ws.Range("F18").Value = 1594.66
$ws.Range("G18").Value = 19.77
$ws.Range("P18").Value = 0.0005935119635064693
$ws.Range("Q18").Value = [double]"2.828718412968693e-08"
$ws.Range("R18").Value = 2.012382810330637
$ws.Range("S18").Value = 3.535149087219765
$ws.Range("T18").Value = 0.02090012160625787
$ws.Range("U18").Value = 2.085471843265076
$ws.Range("V18").Value = [double]"9.117856782426842e-09"
$ws.Range("W18").Value = [double]"7.831861045648451e-05"
$ws.Range("D19").Value = 691195
$ws.Range("E19").Value = 3896072
$ws.Range("F19").Value = 1594.9
$ws.Range("G19").Value = 19.78
$ws.Range("P19").Value = 0.0006913425350652059
$ws.Range("Q19").Value = [double]"3.254298181867987e-08"
$ws.Range("R19").Value = 2.013491670048263
$ws.Range("S19").Value = 3.892948343108125
$ws.Range("T19").Value = 0.02301416172059972
$ws.Range("U19").Value = 2.085484093773392
$ws.Range("V19").Value = [double]"1.05087183821401e-08"
$ws.Range("W19").Value = [double]"8.895157222589057e-05"
$ws.Range("D20").Value = 788963
$ws.Range("E20").Value = 4199748
$ws.Range("F20").Value = 1595.1
$ws.Range("G20").Value = 19.79
$ws.Range("P20").Value = 0.0007891311369653564
$ws.Range("Q20").Value = [double]"3.674993018030371e-08"
$ws.Range("R20").Value = 2.014531483979411
$ws.Range("S20").Value = 4.196546925598039
$ws.Range("T20").Value = 0.02480769349629693
$ws.Range("U20").Value = 2.085492355303519
$ws.Range("V20").Value = [double]"1.188559252519676e-08"
$ws.Range("W20").Value = [double]"9.920694619956655e-05"
$ws.Range("D21").Value = 887174
$ws.Range("E21").Value = 4517998
$ws.Range("F21").Value = 1595.29
$ws.Range("G21").Value = 19.79
$ws.Range("P21").Value = 0.0008873626933204677
$ws.Range("Q21").Value = [double]"4.101326744757732e-08"
$ws.Range("R21").Value = 2.015381791628192
$ws.Range("S21").Value = 4.514552458846584
$ws.Range("T21").Value = 0.02668750172640323
$ws.Range("U21").Value = 2.085494447249077
$ws.Range("V21").Value = [double]"1.327855766630068e-08"
$ws.Range("W21").Value = 0.0001096936629491691
$ws.Range("D22").Value = 984585
$ws.Range("E22").Value = 4826330
$ws.Range("F22").Value = 1595.54
$ws.Range("G22").Value = 19.79
$ws.Range("P22").Value = 0.0009847940822291213
$ws.Range("Q22").Value = [double]"4.524690787590337e-08"
$ws.Range("R22").Value = 2.016109853709128
$ws.Range("S22").Value = 4.822647610700759
$ws.Range("T22").Value = 0.02850875113254024
$ws.Range("U22").Value = 2.085495612408729
$ws.Range("V22").Value = [double]"1.466127880564669e-08"
$ws.Range("W22").Value = 0.0001200896162060546
$ws.Range("D23").Value = -3
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 1594.71
$ws.Range("G23").Value = 19.77
$ws.Range("V23").Value = [double]"1.254146522890935e-09"
$ws.Range("W23").Value = [double]"1.254146522785571e-05"
